$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Date column (BF) holds the game date as literal text, e.g. "2012-04-30".
# Values were previously mangled into "4-30-2011-12" (an artifact of how the
# NBA stats source rendered dates) - correct it back to ISO "YYYY-MM-DD" text.
# NumberFormat is forced to Text first so Excel doesn't reinterpret the
# unambiguous "2012-04-30" string as a date serial, then the style is reset
# back to Normal so no visible/persistent formatting change is introduced.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    $cell.NumberFormat = "@"
    $cell.Value = "2012-04-30"
    $cell.Style = "Normal"
}
